# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout: 基金代码/基金名称/基金规模/
#    股票总仓位/仓位占比/持有市值(亿元)/仓位排名) right after it, rename the copy to
#    "2022-Q1", and overwrite its three data rows with the new quarter's holdings.
# 2. Insert a new top data-row into the "总计" summary sheet for "2022-Q1" and bump
#    the existing running index (column A) down by one for the rows that shifted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet from a copy of "2021-Q4"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newQ = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ.Name = "2022-Q1"

# Row 2: 002666 / 前海开源沪港深创新成长灵活配置混合A
$newQ.Range("B2").Value = "'002666"
$newQ.Range("C2").Value = "前海开源沪港深创新成长灵活配置混合A"
$newQ.Range("D2").Value = "'11.96"
$newQ.Range("E2").Value = "'81.64"
$newQ.Range("F2").Value = "'6.09"
$newQ.Range("G2").Value = "'0.7284"
$newQ.Range("H2").Value = 8

# Row 3: 002667 / 前海开源沪港深创新成长灵活配置混合C
$newQ.Range("B3").Value = "'002667"
$newQ.Range("C3").Value = "前海开源沪港深创新成长灵活配置混合C"
$newQ.Range("D3").Value = "'3.25"
$newQ.Range("E3").Value = "'81.64"
$newQ.Range("F3").Value = "'6.09"
$newQ.Range("G3").Value = "'0.1979"
$newQ.Range("H3").Value = 8

# Row 4: 001037 / 国投瑞银锐意改革灵活配置混合
$newQ.Range("B4").Value = "'001037"
$newQ.Range("C4").Value = "国投瑞银锐意改革灵活配置混合"
$newQ.Range("D4").Value = "'2.60"
$newQ.Range("E4").Value = "'88.37"
$newQ.Range("F4").Value = "'3.48"
$newQ.Range("G4").Value = "'0.0905"
$newQ.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q1 row
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing data rows (2..6) down to (3..7), carrying formatting along.
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

# New first row: index 0, 2022-Q1, 3 funds, 1.02 (亿元)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 1.02

# Bump the running index of the rows that moved down by one.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

Write-Output "done"
